$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells: columns A-J use the "_old" -> "_FV2310" suffix,
#    columns L-U use the "_new" -> "_FV2404" suffix (column K = "diff" is unchanged).
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value().ToString() -replace "_old$", "_FV2310")
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value().ToString() -replace "_new$", "_FV2404")
}

# 2. Convert the data range into a native Excel Table ("Table1") with headers.
#    The header row already carries an explicit bold/fill/border style; stash a
#    copy of that formatting in a scratch range, reset the header to the
#    default style so Excel doesn't synthesize an extra header dxf when the
#    ListObject is created, then restore the original formatting afterwards.
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A100:U100")
$headerRange.Copy($scratchRange) | Out-Null
$headerRange.Style = "Normal"

$dataRange = $ws.Range("A1:U74")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

$scratchRange.Copy($headerRange) | Out-Null
$scratchRange.Clear() | Out-Null

# 3. Freeze the header row (split/freeze after row 1).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "edit complete"
